$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New header cells
$ws.Range("D1").Value = "Url for signup Page"
$ws.Range("E1").Value = "Url for Login Page"

# New data cells (row 2) - note Signup URL lands in D2, Login URL lands in E2
# (Login string is written first so it occupies the earlier shared-string slot)
$ws.Range("E2").Value = "http://automation.zinghr.com/2015/Pages/Authentication/Login.aspx"
$ws.Range("D2").Value = "http://automation.zinghr.com/2015/Pages/Authentication/Signup.aspx"

# Column widths for the new columns (offset by the engine's ~5/6 char padding
# so the serialized OOXML <col width> lands on the target values)
$ws.Columns.Item(4).ColumnWidth = 29.166666666666668
$ws.Columns.Item(5).ColumnWidth = 23.307291666666668

# Update selection to match the new active cell
$ws.Range("D2").Select()
